# Auto update Excel log
# Appends new sensor-log rows to five worksheets, matching the latest
# export from the SeniorConnect logging pipeline (2026-02-01 16:15-16:16).

$wb = $excel.ActiveWorkbook

function Add-LogRow {
    param(
        $ws,
        [int]$Row,
        [string]$Date,
        [string]$Timestamp,
        [string]$Hour,
        [string]$Location,
        $Value,
        [string]$Status
    )

    # Column A holds a date-formatted string ("2026-02-01"); force text so
    # Excel does not silently convert it to a date serial number.
    $dateCell = $ws.Cells.Item($Row, 1)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $Date

    $ws.Cells.Item($Row, 2).Value = $Timestamp
    $ws.Cells.Item($Row, 3).Value = $Hour
    $ws.Cells.Item($Row, 4).Value = $Location
    $ws.Cells.Item($Row, 5).Value = $Value
    $ws.Cells.Item($Row, 6).Value = $Status
}

# --- mmWave(InBed): rows 41-46 --------------------------------------------
$wsInBed = $wb.Worksheets.Item("mmWave(InBed)")
Add-LogRow $wsInBed 41 "2026-02-01" "16:16:21" "16:00" "Bedroom" "Out of Bed" "Empty"
Add-LogRow $wsInBed 42 "2026-02-01" "16:16:22" "16:00" "Bedroom" "In Bed"     "Occupied"
Add-LogRow $wsInBed 43 "2026-02-01" "16:16:26" "16:00" "Bedroom" "In Bed"     "Occupied"
Add-LogRow $wsInBed 44 "2026-02-01" "16:16:26" "16:00" "Bedroom" "In Bed"     "Occupied"
Add-LogRow $wsInBed 45 "2026-02-01" "16:16:27" "16:00" "Bedroom" "In Bed"     "Occupied"
Add-LogRow $wsInBed 46 "2026-02-01" "16:16:28" "16:00" "Bedroom" "In Bed"     "Occupied"

# --- Proximity: rows 22-23 --------------------------------------------------
$wsProximity = $wb.Worksheets.Item("Proximity")
Add-LogRow $wsProximity 22 "2026-02-01" "16:15:50" "16:00" "Living Room Main Door" "ENTER" "User ENTERED Living Room Main Door"
Add-LogRow $wsProximity 23 "2026-02-01" "16:16:16" "16:00" "Bedroom Door"          "ENTER" "User ENTERED Bedroom"

# --- Camera: row 17 ----------------------------------------------------------
$wsCamera = $wb.Worksheets.Item("Camera")
Add-LogRow $wsCamera 17 "2026-02-01" "16:15:52" "16:00" "Living Room Main Door" "Image Captured" "Active"

# --- mmWave(BR): rows 41-46 (numeric Value column) --------------------------
$wsBR = $wb.Worksheets.Item("mmWave(BR)")
Add-LogRow $wsBR 41 "2026-02-01" "16:16:21" "16:00" "Bedroom" 0  "Empty"
Add-LogRow $wsBR 42 "2026-02-01" "16:16:22" "16:00" "Bedroom" 0  "Occupied"
Add-LogRow $wsBR 43 "2026-02-01" "16:16:26" "16:00" "Bedroom" 15 "Occupied"
Add-LogRow $wsBR 44 "2026-02-01" "16:16:27" "16:00" "Bedroom" 2  "Occupied"
Add-LogRow $wsBR 45 "2026-02-01" "16:16:28" "16:00" "Bedroom" 26 "Occupied"
Add-LogRow $wsBR 46 "2026-02-01" "16:16:29" "16:00" "Bedroom" 2  "Occupied"

# --- mmWave(HR): rows 41-46 (numeric Value column) ---------------------------
$wsHR = $wb.Worksheets.Item("mmWave(HR)")
Add-LogRow $wsHR 41 "2026-02-01" "16:16:21" "16:00" "Bedroom" 0  "Empty"
Add-LogRow $wsHR 42 "2026-02-01" "16:16:22" "16:00" "Bedroom" 0  "Occupied"
Add-LogRow $wsHR 43 "2026-02-01" "16:16:26" "16:00" "Bedroom" 63 "Occupied"
Add-LogRow $wsHR 44 "2026-02-01" "16:16:27" "16:00" "Bedroom" 50 "Occupied"
Add-LogRow $wsHR 45 "2026-02-01" "16:16:28" "16:00" "Bedroom" 74 "Occupied"
Add-LogRow $wsHR 46 "2026-02-01" "16:16:29" "16:00" "Bedroom" 50 "Occupied"
